$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($r1, $r2) {
    $a1 = $ws.Cells.Item($r1, 1).Value2
    $b1 = $ws.Cells.Item($r1, 2).Value2
    $c1 = $ws.Cells.Item($r1, 3).Value2

    $a2 = $ws.Cells.Item($r2, 1).Value2
    $b2 = $ws.Cells.Item($r2, 2).Value2
    $c2 = $ws.Cells.Item($r2, 3).Value2

    $ws.Cells.Item($r1, 1).Value2 = $a2
    $ws.Cells.Item($r1, 2).Value2 = $b2
    $ws.Cells.Item($r1, 3).Value2 = $c2

    $ws.Cells.Item($r2, 1).Value2 = $a1
    $ws.Cells.Item($r2, 2).Value2 = $b1
    $ws.Cells.Item($r2, 3).Value2 = $c1
}

Swap-Rows 7 14
Swap-Rows 10 15
Swap-Rows 11 16
